$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 4.854699999999999
$ws.Range("A9").Value = -21.85420000000001
$ws.Range("B12").Value = 5.428799999999993
$ws.Range("C15").Value = -13.57159999999999
$ws.Range("A18").Value = -22.1194
$ws.Range("A20").Value = -19.44229999999999
$ws.Range("B26").Value = 4.056600000000002
$ws.Range("A27").Value = -21.75279999999999
$ws.Range("B27").Value = 5.398800000000006
$ws.Range("B29").Value = 4.921599999999998
$ws.Range("B37").Value = 8.756400000000003
$ws.Range("B38").Value = 4.523500000000001
$ws.Range("C38").Value = -12.51240000000001
$ws.Range("C44").Value = -13.38999999999999
$ws.Range("B51").Value = 5.946000000000002
$ws.Range("C51").Value = -11.8511
$ws.Range("B55").Value = 4.962499999999995
$ws.Range("C57").Value = -13.84009999999999
$ws.Range("C63").Value = -11.7925
$ws.Range("A69").Value = -21.70350000000001
$ws.Range("B69").Value = 5.417999999999994
$ws.Range("B70").Value = 6.078300000000006
$ws.Range("C70").Value = -11.8191
$ws.Range("A76").Value = -19.57629999999998
$ws.Range("A82").Value = -21.9228
$ws.Range("B83").Value = 5.9832
$ws.Range("C99").Value = -12.84899999999999
$ws.Range("B102").Value = 8.191800000000006
